$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 51.93629233333333
$ws.Range("H2").Value = 155.808877
$ws.Range("I2").Value = 0.7704232182162135
$ws.Range("J2").Value = 0.7704232182162134
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 197.1082184852255
$ws.Range("R2").Value = 1773.973966367029
$ws.Range("S2").Value = 0.008044865221818045
$ws.Range("T2").Value = 0.008044865221818043
$ws.Range("G3").Value = 51.93629233333333
$ws.Range("H3").Value = 155.808877
$ws.Range("I3").Value = 0.7704232182162135
$ws.Range("J3").Value = 0.7704232182162134
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 12640.06778818588
$ws.Range("R3").Value = 113760.6100936729
$ws.Range("S3").Value = 0.51589752336086
$ws.Range("T3").Value = 0.5158975233608599
$ws.Range("G4").Value = 51.93629233333333
$ws.Range("H4").Value = 155.808877
$ws.Range("I4").Value = 0.7704232182162135
$ws.Range("J4").Value = 0.7704232182162134
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 1547.772716190122
$ws.Range("R4").Value = 13929.9544457111
$ws.Range("S4").Value = 0.06317150543720274
$ws.Range("T4").Value = 0.06317150543720272
$ws.Range("G5").Value = 51.93629233333333
$ws.Range("H5").Value = 155.808877
$ws.Range("I5").Value = 0.7704232182162135
$ws.Range("J5").Value = 0.7704232182162134
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 4491.283983985057
$ws.Range("R5").Value = 40421.55585586551
$ws.Range("S5").Value = 0.1833093241963327
$ws.Range("T5").Value = 0.1833093241963327
$ws.Range("G6").Value = 0.03171066666666666
$ws.Range("H6").Value = 0.09513199999999999
$ws.Range("I6").Value = 0.0004703961867034368
$ws.Range("J6").Value = 0.0004703961867034368
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 0.1203480790182222
$ws.Range("R6").Value = 1.083132711164
$ws.Range("S6").Value = 0.000004911941687905203
$ws.Range("T6").Value = 0.000004911941687905203
$ws.Range("G7").Value = 0.03171066666666666
$ws.Range("H7").Value = 0.09513199999999999
$ws.Range("I7").Value = 0.0004703961867034368
$ws.Range("J7").Value = 0.0004703961867034368
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("Q7").Value = 7.717627852652445
$ws.Range("S7").Value = 0.0003149908024326838
$ws.Range("T7").Value = 0.0003149908024326838
$ws.Range("G8").Value = 0.03171066666666666
$ws.Range("H8").Value = 0.09513199999999999
$ws.Range("I8").Value = 0.0004703961867034368
$ws.Range("J8").Value = 0.0004703961867034368
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 0.9450213419906666
$ws.Range("R8").Value = 8.505192077915998
$ws.Range("S8").Value = 0.00003857053443272022
$ws.Range("T8").Value = 0.00003857053443272022
$ws.Range("G9").Value = 0.03171066666666666
$ws.Range("H9").Value = 0.09513199999999999
$ws.Range("I9").Value = 0.0004703961867034368
$ws.Range("J9").Value = 0.0004703961867034368
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 2.742236746655111
$ws.Range("R9").Value = 24.680130719896
$ws.Range("S9").Value = 0.0001119229081501275
$ws.Range("T9").Value = 0.0001119229081501275
$ws.Range("G10").Value = 2.162051666666667
$ws.Range("H10").Value = 6.486155
$ws.Range("I10").Value = 0.03207188515291837
$ws.Range("J10").Value = 0.03207188515291837
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 8.205401909603891
$ws.Range("R10").Value = 73.84861718643501
$ws.Range("S10").Value = 0.0003348990364831473
$ws.Range("T10").Value = 0.0003348990364831473
$ws.Range("G11").Value = 2.162051666666667
$ws.Range("H11").Value = 6.486155
$ws.Range("I11").Value = 0.03207188515291837
$ws.Range("J11").Value = 0.03207188515291837
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 526.1923483645978
$ws.Range("R11").Value = 4735.73113528138
$ws.Range("S11").Value = 0.02147625581458147
$ws.Range("T11").Value = 0.02147625581458147
$ws.Range("G12").Value = 2.162051666666667
$ws.Range("H12").Value = 6.486155
$ws.Range("I12").Value = 0.03207188515291837
$ws.Range("J12").Value = 0.03207188515291837
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 64.43210383950166
$ws.Range("R12").Value = 579.888934555515
$ws.Range("S12").Value = 0.002629761434254094
$ws.Range("T12").Value = 0.002629761434254093
$ws.Range("G13").Value = 2.162051666666667
$ws.Range("H13").Value = 6.486155
$ws.Range("I13").Value = 0.03207188515291837
$ws.Range("J13").Value = 0.03207188515291837
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 186.9672937129544
$ws.Range("R13").Value = 1682.70564341659
$ws.Range("S13").Value = 0.007630968867599657
$ws.Range("T13").Value = 0.007630968867599657
$ws.Range("G14").Value = 13.28262333333333
$ws.Range("H14").Value = 39.84787
$ws.Range("I14").Value = 0.1970345004441647
$ws.Range("J14").Value = 0.1970345004441647
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 50.41011024122112
$ws.Range("R14").Value = 453.6909921709901
$ws.Range("S14").Value = 0.002057461357137736
$ws.Range("T14").Value = 0.002057461357137736
$ws.Range("G15").Value = 13.28262333333333
$ws.Range("H15").Value = 39.84787
$ws.Range("I15").Value = 0.1970345004441647
$ws.Range("J15").Value = 0.1970345004441647
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 3232.677031712502
$ws.Range("R15").Value = 29094.09328541252
$ws.Range("S15").Value = 0.1319399628572223
$ws.Range("T15").Value = 0.1319399628572223
$ws.Range("G16").Value = 13.28262333333333
$ws.Range("H16").Value = 39.84787
$ws.Range("I16").Value = 0.1970345004441647
$ws.Range("J16").Value = 0.1970345004441647
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 395.8403858099233
$ws.Range("R16").Value = 3562.56347228931
$ws.Range("S16").Value = 0.01615601103630282
$ws.Range("T16").Value = 0.01615601103630281
$ws.Range("G17").Value = 13.28262333333333
$ws.Range("H17").Value = 39.84787
$ws.Range("I17").Value = 0.1970345004441647
$ws.Range("J17").Value = 0.1970345004441647
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 1148.638664066096
$ws.Range("R17").Value = 10337.74797659486
$ws.Range("S17").Value = 0.04688106519350191
$ws.Range("T17").Value = 0.04688106519350191
